$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: ref -> new text value.
# A leading apostrophe forces Excel to store the value as literal text
# (otherwise numeric-looking strings like "331.60" or "0.49%" get parsed
# into numbers), and resetting the Style back to "Normal" afterwards avoids
# leaving a stray quote-prefixed number format/style on the cell.
$updates = @(
    @{Cell='D2'; Value='331.60'}
    @{Cell='E2'; Value='0.49%'}
    @{Cell='D3'; Value='45.51'}
    @{Cell='E3'; Value='3.20%'}
    @{Cell='D4'; Value='5.611'}
    @{Cell='E4'; Value='2.12%'}
    @{Cell='D5'; Value='0.08349'}
    @{Cell='E5'; Value='4.30%'}
    @{Cell='D6'; Value='2.100'}
    @{Cell='E6'; Value='6.33%'}
    @{Cell='D7'; Value='0.9625'}
    @{Cell='E7'; Value='1.12%'}
    @{Cell='D8'; Value='2.555'}
    @{Cell='E8'; Value='-0.73%'}
    @{Cell='D9'; Value='0.1160'}
    @{Cell='E9'; Value='3.88%'}
    @{Cell='D10'; Value='0.1928'}
    @{Cell='E10'; Value='2.02%'}
    @{Cell='D11'; Value='10.38'}
    @{Cell='E11'; Value='-1.89%'}
    @{Cell='D12'; Value='0.09848'}
    @{Cell='E12'; Value='-0.66%'}
    @{Cell='D13'; Value='0.04611'}
    @{Cell='E13'; Value='-3.65%'}
    @{Cell='D14'; Value='0.1060'}
    @{Cell='E14'; Value='-0.50%'}
    @{Cell='D15'; Value='0.001297'}
    @{Cell='E15'; Value='2.40%'}
    @{Cell='D16'; Value='0.006056'}
    @{Cell='E16'; Value='1.01%'}
    @{Cell='B17'; Value='LEO'}
    @{Cell='C17'; Value='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'}
    @{Cell='D17'; Value='3.376'}
    @{Cell='E17'; Value='0.20%'}
    @{Cell='B18'; Value='GateToken'}
    @{Cell='C18'; Value='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'}
    @{Cell='D18'; Value='4.435'}
    @{Cell='E18'; Value='1.05%'}
    @{Cell='B19'; Value='BitpandaEcosystemToken'}
    @{Cell='C19'; Value='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'}
    @{Cell='D19'; Value='0.3343'}
    @{Cell='E19'; Value='-3.75%'}
    @{Cell='B20'; Value='ProBitToken'}
    @{Cell='C20'; Value='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'}
    @{Cell='D20'; Value='0.1392'}
    @{Cell='E20'; Value='-1.67%'}
    @{Cell='B21'; Value='ZBToken'}
    @{Cell='C21'; Value='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'}
    @{Cell='D21'; Value='0.2653'}
    @{Cell='E21'; Value='2.54%'}
    @{Cell='B22'; Value='CoinExToken'}
    @{Cell='C22'; Value='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'}
    @{Cell='D22'; Value='0.04173'}
    @{Cell='E22'; Value='2.16%'}
    @{Cell='B23'; Value='BitKan'}
    @{Cell='C23'; Value='https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'}
    @{Cell='D23'; Value='0.001318'}
    @{Cell='E23'; Value='3.50%'}
    @{Cell='B24'; Value='HotbitToken'}
    @{Cell='C24'; Value='https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'}
    @{Cell='D24'; Value='0.004559'}
    @{Cell='E24'; Value='4.20%'}
    @{Cell='E25'; Value='8.50%'}
    @{Cell='D26'; Value='0.0003748'}
    @{Cell='E26'; Value='0.06%'}
    @{Cell='D38'; Value='0.02699'}
    @{Cell='E38'; Value='4.19%'}
    @{Cell='D39'; Value='0.05757'}
    @{Cell='E39'; Value='1.01%'}
    @{Cell='D40'; Value='0.007844'}
    @{Cell='E40'; Value='3.76%'}
    @{Cell='D41'; Value='0.1433'}
    @{Cell='E41'; Value='2.35%'}
    @{Cell='D42'; Value='0.007247'}
    @{Cell='E42'; Value='-1.47%'}
    @{Cell='D43'; Value='0.002024'}
    @{Cell='E43'; Value='0.41%'}
    @{Cell='D44'; Value='0.008873'}
    @{Cell='E44'; Value='6.15%'}
    @{Cell='D45'; Value='0.3544'}
    @{Cell='D46'; Value='0.00007109'}
    @{Cell='E46'; Value='-0.36%'}
    @{Cell='E47'; Value='0.16%'}
    @{Cell='D48'; Value='0.0005813'}
    @{Cell='E48'; Value='0.23%'}
    @{Cell='D49'; Value='0.003505'}
    @{Cell='E49'; Value='-1.32%'}
    @{Cell='D50'; Value='0.003506'}
    @{Cell='E50'; Value='-0.73%'}
    @{Cell='D51'; Value='0.00002104'}
    @{Cell='E51'; Value='0.16%'}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
